# Update the "108th IETF Online" footers to "109th IETF Online" on every
# slide, and add a new bullet ("Need for different HbH and E2E Indicator
# Labels") just above "Added MSD consideration" on the
# "Updates Since IETF-108 (Version-02)" slide.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)

        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }

        $tr = $sh.TextFrame.TextRange

        # Footer placeholders hold a run with exactly "108" immediately
        # followed by a superscript "th" run - bump just that run's text
        # so the "th"/" IETF Online" runs (and their formatting) survive.
        if ($tr.Text -eq "108th IETF Online") {
            $tr.Characters(1, 3).Text = "109"
        }
    }
}

# Slide 5 ("Updates Since IETF-108 (Version-02)") content placeholder:
# insert a new sub-bullet above "Added MSD consideration".
$slide5 = $p.Slides.Item(5)
for ($j = 1; $j -le $slide5.Shapes.Count; $j++) {
    $sh = $slide5.Shapes.Item($j)
    if ($sh.Name -like "Content Placeholder*") {
        $tr = $sh.TextFrame.TextRange
        for ($k = 1; $k -le $tr.Paragraphs().Count; $k++) {
            $para = $tr.Paragraphs($k)
            if ($para.Text.TrimEnd("`r") -eq "Added MSD consideration") {
                $newParaStart = $para.Start

                # Insert the new bullet's full text (plus trailing CR to
                # create a new paragraph) ahead of "Added MSD consideration".
                # It inherits that paragraph's formatting (lvl 1, Wingdings
                # "ü" bullet).
                [void]$para.InsertBefore("Need for different HbH and E2E Indicator Labels`r")

                $full = $sh.TextFrame.TextRange

                $part1 = "Need for different "
                $part2 = "HbH"
                $part3 = " and E2E Indicator Labels"

                # Re-split the inserted text into three runs so "HbH" can
                # carry its own run (matching the source formatting split).
                # Re-assigning a sub-range's Text to its own value forces the
                # engine to materialize it as a distinct run without
                # introducing any incidental formatting attributes.
                $r2 = $full.Characters($newParaStart + $part1.Length, $part2.Length)
                $r2.Text = $r2.Text

                $r3 = $full.Characters($newParaStart + $part1.Length + $part2.Length, $part3.Length)
                $r3.Text = $r3.Text

                break
            }
        }
        break
    }
}
